# Insert a new data row at row 154, shifting the existing rows 154:198 down to
# 155:199, then populate the newly inserted row with the new record
# (fecha 44988 / Región Metropolitana, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 154, pushing everything below it down
# by one (old row 154 becomes 155, ..., old row 198 becomes 199).
$ws.Rows.Item(154).Insert()

# Populate the new row 154 with the new record's data.
$ws.Range("A154").Value = 11
$ws.Range("B154").Value = "Vega Monumental Concepción"
$ws.Range("C154").Value = "Bíobío"
$ws.Range("D154").Value = 44988
$ws.Range("E154").Value = 8
$ws.Range("F154").Value = 100112032
$ws.Range("G154").Value = "Zapallo italiano"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 220
$ws.Range("K154").Value = 7000
$ws.Range("L154").Value = 8000
$ws.Range("M154").Value = 7455
$ws.Range("N154").Value = "$/caja 60 unidades"
$ws.Range("O154").Value = "Región Metropolitana"
$ws.Range("P154").Value = 124
$ws.Range("Q154").Value = 60
$ws.Range("R154").Value = "Hortaliza"
